$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2615.0625
$ws.Range("I38").Value = 125.46667
$ws.Range("J38").Value = 4811.7646
$ws.Range("K38").Value = 376.40001
$ws.Range("L38").Value = 14435.2938
$ws.Range("M38").Value = -4.400009999999952
$ws.Range("N38").Value = -15179.2938

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 28266.445
$ws.Range("J64").Value = 28266.445
$ws.Range("L64").Value = 28266.445
$ws.Range("N64").Value = -28762.445

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 28266.445
$ws.Range("J67").Value = 28266.445
$ws.Range("L67").Value = 28266.445
$ws.Range("N67").Value = -29982.445

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 34023.625
$ws.Range("I135").Value = 1151.1
$ws.Range("J135").Value = 88811.164
$ws.Range("K135").Value = 10359.9
$ws.Range("L135").Value = 799300.476
$ws.Range("M135").Value = -7824.9
$ws.Range("N135").Value = -804370.476

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5171.5
$ws.Range("I137").Value = 2560.923
$ws.Range("J137").Value = 8942.333000000001
$ws.Range("K137").Value = 7682.768999999999
$ws.Range("L137").Value = 26826.999
$ws.Range("M137").Value = -5132.768999999999
$ws.Range("N137").Value = -31926.999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4474.4844
$ws.Range("J138").Value = 3576.3726
$ws.Range("L138").Value = 10729.1178
$ws.Range("N138").Value = -21009.1178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 717618.0600000001
$ws.Range("I45").Value = 1668496.5
$ws.Range("J45").Value = 4459.25
$ws.Range("K45").Value = 1668496.5
$ws.Range("L45").Value = 4459.25
$ws.Range("M45").Value = -1668119.5
$ws.Range("N45").Value = -5213.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3356.4146
$ws.Range("I134").Value = 3239.842
$ws.Range("K134").Value = 9719.526
$ws.Range("M134").Value = -7184.526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 50863.19
$ws.Range("I134").Value = 58876.555
$ws.Range("K134").Value = 176629.665
$ws.Range("M134").Value = -174094.665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 10999
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 10999
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 893.8570999999999
$ws.Range("J68").Value = 815
$ws.Range("L68").Value = 2445
$ws.Range("N68").Value = -4067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 893.8570999999999
$ws.Range("J71").Value = 815
$ws.Range("L71").Value = 7335
$ws.Range("N71").Value = -15447

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 12000
$ws.Range("J74").Value = 12000
$ws.Range("L74").Value = 36000
$ws.Range("N74").Value = -38122

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1854.6666
$ws.Range("J75").Value = 2025.6
$ws.Range("L75").Value = 6076.799999999999
$ws.Range("N75").Value = -8072.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 12000
$ws.Range("J77").Value = 12000
$ws.Range("L77").Value = 108000
$ws.Range("N77").Value = -118608

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 1854.6666
$ws.Range("J78").Value = 2025.6
$ws.Range("L78").Value = 18230.4
$ws.Range("N78").Value = -28214.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3763.6924
$ws.Range("I80").Value = 3866
$ws.Range("J80").Value = 3745.0908
$ws.Range("K80").Value = 11598
$ws.Range("L80").Value = 11235.2724
$ws.Range("M80").Value = -10662
$ws.Range("N80").Value = -13107.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 3763.6924
$ws.Range("I83").Value = 3866
$ws.Range("J83").Value = 3745.0908
$ws.Range("K83").Value = 34794
$ws.Range("L83").Value = 33705.8172
$ws.Range("M83").Value = -30114
$ws.Range("N83").Value = -43065.8172

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 152603.84
$ws.Range("I140").Value = 160477.73
$ws.Range("K140").Value = 481433.1900000001
$ws.Range("M140").Value = -476253.1900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2786.6667
$ws.Range("I80").Value = 3243.6
$ws.Range("J80").Value = 2460.2856
$ws.Range("K80").Value = 3243.6
$ws.Range("L80").Value = 2460.2856
$ws.Range("M80").Value = -2245.6
$ws.Range("N80").Value = -4456.2856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2786.6667
$ws.Range("I83").Value = 3243.6
$ws.Range("J83").Value = 2460.2856
$ws.Range("K83").Value = 16218
$ws.Range("L83").Value = 12301.428
$ws.Range("M83").Value = -11226
$ws.Range("N83").Value = -22285.428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4266.606
$ws.Range("I122").Value = 2753.8462
$ws.Range("J122").Value = 5249.9
$ws.Range("K122").Value = 8261.5386
$ws.Range("L122").Value = 15749.7
$ws.Range("M122").Value = -5811.5386
$ws.Range("N122").Value = -20649.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 567534.7
$ws.Range("I7").Value = 637651.5
$ws.Range("J7").Value = 6600
$ws.Range("K7").Value = 637651.5
$ws.Range("L7").Value = 6600
$ws.Range("M7").Value = -637539.5
$ws.Range("N7").Value = -6824

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3117.6875
$ws.Range("I22").Value = 1289.8572
$ws.Range("J22").Value = 4539.3335
$ws.Range("K22").Value = 1289.8572
$ws.Range("L22").Value = 4539.3335
$ws.Range("M22").Value = -994.8571999999999
$ws.Range("N22").Value = -5129.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3117.6875
$ws.Range("I27").Value = 1289.8572
$ws.Range("J27").Value = 4539.3335
$ws.Range("K27").Value = 1289.8572
$ws.Range("L27").Value = 4539.3335
$ws.Range("M27").Value = -1182.8572
$ws.Range("N27").Value = -4753.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1111
$ws.Range("I55").Value = 331.4
$ws.Range("J55").Value = 2085.5
$ws.Range("K55").Value = 331.4
$ws.Range("L55").Value = 2085.5
$ws.Range("M55").Value = -158.4
$ws.Range("N55").Value = -2431.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4379.775
$ws.Range("I122").Value = 3818.6667
$ws.Range("K122").Value = 11456.0001
$ws.Range("M122").Value = -9006.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 567534.7
$ws.Range("I126").Value = 637651.5
$ws.Range("J126").Value = 6600
$ws.Range("K126").Value = 1912954.5
$ws.Range("L126").Value = 19800
$ws.Range("M126").Value = -1910484.5
$ws.Range("N126").Value = -24740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 42462.84
$ws.Range("I126").Value = 54904.58
$ws.Range("J126").Value = 3064
$ws.Range("K126").Value = 164713.74
$ws.Range("L126").Value = 9192
$ws.Range("M126").Value = -162243.74
$ws.Range("N126").Value = -14132

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2809.851
$ws.Range("I136").Value = 2411.0476
$ws.Range("J136").Value = 6159.8
$ws.Range("K136").Value = 7233.1428
$ws.Range("L136").Value = 18479.4
$ws.Range("M136").Value = -4683.1428
$ws.Range("N136").Value = -23579.4
